$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.88
$wsSummary.Range("B4").Value = -0.12
$wsSummary.Range("B5").Value = -0.24
$wsSummary.Range("B6").Value = 10
$wsSummary.Range("B8").Value = 7
$wsSummary.Range("B9").Value = 30

# --- Strategy Status sheet (MarketMaking row) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.88
$wsStatus.Range("D4").Value = 10
$wsStatus.Range("E4").Value = -0.12
$wsStatus.Range("F4").Value = -0.12
$wsStatus.Range("G4").Value = 30

# --- New trade row (Trade #10), appended to both "All Trades" and "MarketMaking" sheets ---
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A11").Value = 10

    # Force the date column to stay plain text (matches the source data,
    # which stores dates as literal strings, not Excel date serials).
    $ws.Range("B11").NumberFormat = "@"
    $ws.Range("B11").Value = "2026-02-17"

    $ws.Range("C11").Value = "13:34:23"
    $ws.Range("D11").Value = "MarketMaking"
    $ws.Range("E11").Value = "UP"
    $ws.Range("F11").Value = 0.08
    $ws.Range("G11").Value = 0.040657
    $ws.Range("H11").Value = "CLOSED"
    $ws.Range("I11").Value = -49.1787
    $ws.Range("J11").Value = -0.04
    $ws.Range("K11").Value = 99.88
    $ws.Range("L11").Value = 0
    $ws.Range("M11").Value = 0
    $ws.Range("N11").Value = 0.6
    $ws.Range("O11").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P11").Value = "early_exit"
    $ws.Range("Q11").Value = 0.13
}
